$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '59.412.49'
$ws.Range("E2").Value = '  -2.65%  '
$ws.Range("D3").Value = '2.362.26'
$ws.Range("E3").Value = '  -2.86%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '558.03'
$ws.Range("E5").Value = '  -2.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.30'
$ws.Range("E6").Value = '  -2.27%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.527'
$ws.Range("E8").Value = '  -0.50%  '
$ws.Range("D9").Value = '2.356.40'
$ws.Range("E9").Value = '  -2.55%  '
$ws.Range("E10").Value = '  -3.96%  '
$ws.Range("E11").Value = '  -1.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.05'
$ws.Range("E12").Value = '  -1.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.335'
$ws.Range("E13").Value = '  -1.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.39'
$ws.Range("E14").Value = '  -2.96%  '
$ws.Range("E16").Value = '  -3.81%  '
$ws.Range("D17").Value = '59.286.60'
$ws.Range("E17").Value = '  -2.86%  '
$ws.Range("D18").Value = '2.336.13'
$ws.Range("E18").Value = '  -4.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.92'
$ws.Range("E19").Value = '  +8.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.40'
$ws.Range("E20").Value = '  -1.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '320.50'
$ws.Range("E21").Value = '  -1.04%  '
$ws.Range("E22").Value = '  -0.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.93'
$ws.Range("E23").Value = '  -3.39%  '
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("E25").Value = '  -5.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '63.97'
$ws.Range("E26").Value = '  -1.57%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.10'
$ws.Range("E27").Value = '  -9.15%  '
$ws.Range("B28").Value = 'Bittensor'
$ws.Range("C28").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '549.70'
$ws.Range("E28").Value = '  -4.04%  '
$ws.Range("E29").Value = '  -3.30%  '
$ws.Range("D30").Value = '0.0₃0910'
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.92'
$ws.Range("E31").Value = '  +0.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.29'
$ws.Range("E32").Value = '  -3.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.76'
$ws.Range("E33").Value = '  -4.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.129'
$ws.Range("E34").Value = '  -2.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.40'
$ws.Range("E36").Value = '  +1.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '149.85'
$ws.Range("E37").Value = '  -0.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.364'
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.51'
$ws.Range("E39").Value = '  -2.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.99'
$ws.Range("E40").Value = '  -1.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.98'
$ws.Range("E41").Value = '  -2.59%  '
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.41'
$ws.Range("E43").Value = '  -0.69%  '
$ws.Range("E44").Value = '  -1.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.39'
$ws.Range("E45").Value = '  +1.90%  '
$ws.Range("D46").Value = '0.0₆0288'
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '137.65'
$ws.Range("E47").Value = '  -2.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.48'
$ws.Range("E48").Value = '  -1.23%  '
$ws.Range("E49").Value = '  -1.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0497'
$ws.Range("E50").Value = '  -1.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.98'
$ws.Range("E51").Value = '  -2.91%  '
